$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3 with new record data (replaces previous row 3 contents) ---
$ws.Range("A3").Value = "'3765"
$ws.Range("B3").Value = "'10/16/2024"
$ws.Range("C3").Value = "'NAZCA AV. 1675"
$ws.Range("E3").Value = "'01082878"
$ws.Range("H3").Value = "'ya se traspaso nodo retirar columna"
$ws.Range("K3").Value = "'Fuente Teco"
$ws.Range("M3").Value = -58.47874
$ws.Range("N3").Value = -34.61462
$ws.Range("Q3").Value = "'NRA-M"

# --- Remove the old row 49 record (SAN NICOLAS 5045 / Caso 7169) ---
# This shifts every following row up by one, which matches the new
# ordering for rows 49-68, and updates the sheet dimension to A1:R68.
$ws.Rows(49).Delete()
